# Insert one new data row at row 87 (weekly "Repollo" price record for
# Macroferia Regional de Talca), pushing the existing rows 87-180 down to
# 88-181 and extending the sheet's used range to A1:R181.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 87; Excel shifts rows 87..180
# down to 88..181 automatically (mirrors the row-insert behaviour implied
# by the diff, where every row from 87 onward simply moves down by one).
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record.
$ws.Range("A87").Value = 5
$ws.Range("B87").Value = 'Macroferia Regional de Talca'
$ws.Range("C87").Value = 'Maule'
$ws.Range("D87").Value = 44484
$ws.Range("E87").Value = 7
$ws.Range("F87").Value = 100112006
$ws.Range("G87").Value = 'Repollo'
$ws.Range("H87").Value = 'Crespo record'
$ws.Range("I87").Value = 'Primera'
$ws.Range("J87").Value = 3000
$ws.Range("K87").Value = 700
$ws.Range("L87").Value = 700
$ws.Range("M87").Value = 700
$ws.Range("N87").Value = '$/unidad'
$ws.Range("O87").Value = 'Región del Maule'
$ws.Range("P87").Value = 700
$ws.Range("Q87").Value = 1
$ws.Range("R87").Value = 'Hortaliza'
